# The "Apple" row's price (D3) changes from the number 345 to the literal
# text "600". Prefixing with an apostrophe forces Excel to store it as a
# text value (shared string) rather than re-parsing "600" back into a
# number. Resetting the style to "Normal" afterwards clears the
# quote-prefix formatting flag that Excel would otherwise stamp onto the
# cell, leaving it on the same default style as before the edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D3").Value = "'600"
$ws.Range("D3").Style = "Normal"
